# cyclone-iii-pinout.xlsx: mark several previously-unused pins as
# "Switch Input" / "Extra Switch Input" connections (Active=H, Direction=I).
# Rows 64, 85, 136-138, 184, 185 (column C/D/E/F) go from the blank "-"
# placeholder to H / I / <connection text> / <use text>.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-PinRow {
    param([int]$Row, [string]$Connection)
    $ws.Range("C$Row").Value = "H"
    $ws.Range("D$Row").Value = "I"
    $ws.Range("E$Row").Value = $Connection
    $ws.Range("F$Row").Value = $Connection
}

Set-PinRow 64  "Switch Input"
Set-PinRow 85  "Extra Switch Input"
Set-PinRow 136 "Extra Switch Input"
Set-PinRow 137 "Extra Switch Input"
Set-PinRow 138 "Extra Switch Input"
Set-PinRow 184 "Extra Switch Input"
Set-PinRow 185 "Extra Switch Input"

# Reflect the scrolled/selected cell recorded in the workbook view.
$ws.Range("L173").Select()
